$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.430.28'
$ws.Range('E2').Value = '  +0.88%  '

$ws.Range('D3').Value = '3.367.35'
$ws.Range('E3').Value = '  +0.31%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = "'590.23"
$ws.Range('E5').Value = '  +4.95%  '

$ws.Range('D6').Value = "'187.77"
$ws.Range('E6').Value = '  -1.60%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = "'0.597"
$ws.Range('E8').Value = '  +2.08%  '

$ws.Range('E9').Value = '  -0.58%  '

$ws.Range('D10').Value = "'0.586"
$ws.Range('E10').Value = '  -0.31%  '

$ws.Range('D11').Value = "'47.43"
$ws.Range('E11').Value = '  +0.40%  '

$ws.Range('E12').Value = '  +0.70%  '

$ws.Range('D13').Value = '3.911.48'
$ws.Range('E13').Value = '  +0.52%  '

$ws.Range('D14').Value = "'638.93"
$ws.Range('E14').Value = '  +5.29%  '

$ws.Range('D15').Value = "'8.61"
$ws.Range('E15').Value = '  -1.48%  '

$ws.Range('D16').Value = '67.462.15'
$ws.Range('E16').Value = '  +1.20%  '

$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = "'0.119"

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.372.03'
$ws.Range('E18').Value = '  +0.64%  '

$ws.Range('D19').Value = "'18.00"
$ws.Range('E19').Value = '  -0.49%  '

$ws.Range('E20').Value = '  +0.46%  '

$ws.Range('D21').Value = "'0.909"
$ws.Range('E21').Value = '  +0.02%  '

$ws.Range('D22').Value = "'18.06"
$ws.Range('E22').Value = '  -2.64%  '

$ws.Range('E23').Value = '  +0.55%  '

$ws.Range('D24').Value = "'99.36"
$ws.Range('E24').Value = '  -1.52%  '

$ws.Range('E25').Value = '  -0.42%  '

$ws.Range('E26').Value = '  +2.47%  '

$ws.Range('D27').Value = "'9.72"
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').Value = "'32.41"
$ws.Range('E28').Value = '  +5.14%  '

$ws.Range('D29').Value = "'8.70"
$ws.Range('E29').Value = '  -0.93%  '

$ws.Range('D30').Value = "'6.89"
$ws.Range('E30').Value = '  +1.01%  '

$ws.Range('D31').Value = "'612.69"
$ws.Range('E31').Value = '  +4.35%  '

$ws.Range('D32').Value = "'3.86"
$ws.Range('E32').Value = '  -3.13%  '

$ws.Range('D33').Value = "'11.11"
$ws.Range('E33').Value = '  +0.09%  '

$ws.Range('D34').Value = '3.911.95'
$ws.Range('E34').Value = '  +5.23%  '

$ws.Range('D35').Value = "'0.107"
$ws.Range('E35').Value = '  +1.06%  '

$ws.Range('D36').Value = "'0.998"
$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('D37').Value = "'55.94"
$ws.Range('E37').Value = '  -2.28%  '

$ws.Range('D38').Value = "'2.84"
$ws.Range('E38').Value = '  +4.40%  '

$ws.Range('E39').Value = '  +1.99%  '

$ws.Range('D40').Value = "'3.29"
$ws.Range('E40').Value = '  +0.17%  '

$ws.Range('D41').Value = "'33.71"
$ws.Range('E41').Value = '  -0.93%  '

$ws.Range('E42').Value = '  -1.49%  '

$ws.Range('E43').Value = '  +0.13%  '

$ws.Range('D44').Value = "'3.40"
$ws.Range('E44').Value = '  +0.20%  '

$ws.Range('D45').Value = "'0.0423"
$ws.Range('E45').Value = '  -0.31%  '

$ws.Range('E46').Value = '  -0.38%  '

$ws.Range('E47').Value = '  -0.82%  '

$ws.Range('E48').Value = '  +0.65%  '

$ws.Range('D49').Value = "'1.36"
$ws.Range('E49').Value = '  +9.41%  '

$ws.Range('D50').Value = "'2.83"
$ws.Range('E50').Value = '  -21.98%  '

$ws.Range('D51').Value = "'129.81"
$ws.Range('E51').Value = '  +4.23%  '

